# Weekly update to the Arándano (blue) / Macroferia Regional de Talca dataset.
# Two new rows of observations are inserted at the top of the data block
# (row 50), pushing the existing rows down by two positions, and the new
# rows are populated with this week's figures.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows right before the current row 50; everything that
# used to live at row 50 onward shifts down to row 52 onward, which is
# exactly the data movement seen in the diff (old row N -> new row N+2).
$ws.Rows("50:51").Insert()

# --- New row 50 ---------------------------------------------------------
$ws.Range("A50").Value = 5
$ws.Range("B50").Value = "Macroferia Regional de Talca"
$ws.Range("C50").Value = "Maule"
$ws.Range("D50").Value = 44967
$ws.Range("E50").Value = 7
$ws.Range("F50").Value = "Fruta"
$ws.Range("G50").Value = 100101
$ws.Range("H50").Value = "Berries"
$ws.Range("I50").Value = 100101001
$ws.Range("J50").Value = "Arándano (blue)"
$ws.Range("K50").Value = "Sin especificar"
$ws.Range("L50").Value = "Primera"
$ws.Range("M50").Value = 200
$ws.Range("N50").Value = 3000
$ws.Range("O50").Value = 3000
$ws.Range("P50").Value = 3000
$ws.Range("Q50").Value = "$/bandeja 2 kilos"
$ws.Range("R50").Value = "Provincia de Curicó"
$ws.Range("S50").Value = 1500
$ws.Range("T50").Value = 2

# --- New row 51 ---------------------------------------------------------
$ws.Range("A51").Value = 5
$ws.Range("B51").Value = "Macroferia Regional de Talca"
$ws.Range("C51").Value = "Maule"
$ws.Range("D51").Value = 44967
$ws.Range("E51").Value = 7
$ws.Range("F51").Value = "Fruta"
$ws.Range("G51").Value = 100101
$ws.Range("H51").Value = "Berries"
$ws.Range("I51").Value = 100101001
$ws.Range("J51").Value = "Arándano (blue)"
$ws.Range("K51").Value = "Sin especificar"
$ws.Range("L51").Value = "Segunda"
$ws.Range("M51").Value = 150
$ws.Range("N51").Value = 2500
$ws.Range("O51").Value = 2500
$ws.Range("P51").Value = 2500
$ws.Range("Q51").Value = "$/bandeja 2 kilos"
$ws.Range("R51").Value = "Provincia de Curicó"
$ws.Range("S51").Value = 1250
$ws.Range("T51").Value = 2
